# Refresh the Price (D) and Volume(1h) (E) columns with the latest scraped figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.183.98"
$ws.Range("E2").Value = "  -2.46%  "
$ws.Range("D3").Value = "2.870.35"
$ws.Range("E3").Value = "  -2.32%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'565.46"
$ws.Range("E5").Value = "  -4.62%  "
$ws.Range("D6").Value = "'142.69"
$ws.Range("E6").Value = "  -2.88%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").Value = "'0.504"
$ws.Range("E8").Value = "  -0.68%  "
$ws.Range("D9").Value = "2.870.87"
$ws.Range("E9").Value = "  -2.33%  "
$ws.Range("D10").Value = "'6.83"
$ws.Range("E10").Value = "  -6.48%  "
$ws.Range("E11").Value = "  -3.26%  "
$ws.Range("D12").Value = "'0.430"
$ws.Range("E12").Value = "  -2.58%  "
$ws.Range("D13").Value = "'0.0000233"
$ws.Range("E13").Value = "  -1.68%  "
$ws.Range("D14").Value = "'31.74"
$ws.Range("E14").Value = "  -3.07%  "
$ws.Range("E15").Value = "  -0.59%  "
$ws.Range("D16").Value = "3.347.89"
$ws.Range("E16").Value = "  -2.23%  "
$ws.Range("D17").Value = "61.263.45"
$ws.Range("E17").Value = "  -2.24%  "
$ws.Range("D18").Value = "'6.55"
$ws.Range("E18").Value = "  -2.00%  "
$ws.Range("D19").Value = "2.866.28"
$ws.Range("E19").Value = "  -3.07%  "
$ws.Range("D20").Value = "'430.68"
$ws.Range("E20").Value = "  -2.31%  "
$ws.Range("D21").Value = "'13.02"
$ws.Range("E21").Value = "  -2.79%  "
$ws.Range("D22").Value = "'0.652"
$ws.Range("E22").Value = "  -1.94%  "
$ws.Range("D23").Value = "'6.80"
$ws.Range("E23").Value = "  -3.16%  "
$ws.Range("D24").Value = "'78.85"
$ws.Range("E24").Value = "  -2.95%  "
$ws.Range("D25").Value = "'11.68"
$ws.Range("E25").Value = "  -0.33%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").Value = "'9.96"
$ws.Range("E27").Value = "  -10.25%  "
$ws.Range("D28").Value = "'2.00"
$ws.Range("E28").Value = "  -6.38%  "
$ws.Range("D29").Value = "'0.0000104"
$ws.Range("E29").Value = "  +2.33%  "
$ws.Range("D30").Value = "'6.98"
$ws.Range("E30").Value = "  -3.16%  "
$ws.Range("D31").Value = "'2.49"
$ws.Range("E31").Value = "  -4.54%  "
$ws.Range("E32").Value = "  -8.60%  "
$ws.Range("D33").Value = "'0.999"
$ws.Range("E33").Value = "  -0.21%  "
$ws.Range("E34").Value = "  -2.74%  "
$ws.Range("D35").Value = "'25.42"
$ws.Range("E35").Value = "  -3.54%  "
$ws.Range("D36").Value = "'0.951"
$ws.Range("E36").Value = "  -3.91%  "
$ws.Range("D37").Value = "'5.36"
$ws.Range("E37").Value = "  -3.92%  "
$ws.Range("D38").Value = "'48.78"
$ws.Range("E38").Value = "  -1.73%  "
$ws.Range("E39").Value = "  -4.80%  "
$ws.Range("E40").Value = "  -10.24%  "
$ws.Range("D41").Value = "'8.20"
$ws.Range("E41").Value = "  -3.15%  "
$ws.Range("E42").Value = "  -3.45%  "
$ws.Range("D43").Value = "'39.19"
$ws.Range("E43").Value = "  -0.15%  "
$ws.Range("E44").Value = "  -4.85%  "
$ws.Range("D45").Value = "2.681.00"
$ws.Range("E45").Value = "  -0.72%  "
$ws.Range("D46").Value = "'133.52"
$ws.Range("E46").Value = "  -0.96%  "
$ws.Range("E47").Value = "  -1.25%  "
$ws.Range("D49").Value = "'336.92"
$ws.Range("E49").Value = "  -6.93%  "
$ws.Range("E50").Value = "  -1.90%  "
$ws.Range("D51").Value = "'21.41"
$ws.Range("E51").Value = "  -6.11%  "
